$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 484.9524
$ws.Range("I33").Value = 257
$ws.Range("K33").Value = 257
$ws.Range("M33").Value = -28
$ws.Range("H74").Value = 7791.375
$ws.Range("I74").Value = 10216.2
$ws.Range("K74").Value = 10216.2
$ws.Range("M74").Value = -9280.200000000001
$ws.Range("H77").Value = 7791.375
$ws.Range("I77").Value = 10216.2
$ws.Range("K77").Value = 51081
$ws.Range("M77").Value = -46401
$ws.Range("H129").Value = 3618.5
$ws.Range("I129").Value = 0
$ws.Range("K129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("H137").Value = 1587.2667
$ws.Range("I137").Value = 1450
$ws.Range("J137").Value = 1744.1428
$ws.Range("K137").Value = 4350
$ws.Range("L137").Value = 5232.428400000001
$ws.Range("M137").Value = -1800
$ws.Range("N137").Value = -10332.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20738654
$ws.Range("I32").Value = 19919804
$ws.Range("K32").Value = 19919804
$ws.Range("M32").Value = -19919517
$ws.Range("H61").Value = 3396.742
$ws.Range("I61").Value = 2967.2
$ws.Range("K61").Value = 2967.2
$ws.Range("M61").Value = -2755.2
$ws.Range("H63").Value = 5236.364
$ws.Range("I63").Value = 2650
$ws.Range("J63").Value = 6714.2856
$ws.Range("K63").Value = 2650
$ws.Range("L63").Value = 6714.2856
$ws.Range("M63").Value = -1964
$ws.Range("N63").Value = -8086.2856
$ws.Range("H66").Value = 5236.364
$ws.Range("I66").Value = 2650
$ws.Range("J66").Value = 6714.2856
$ws.Range("K66").Value = 13250
$ws.Range("L66").Value = 33571.428
$ws.Range("M66").Value = -9818
$ws.Range("N66").Value = -40435.428
$ws.Range("H74").Value = 1476.9714
$ws.Range("J74").Value = 2022.5
$ws.Range("L74").Value = 2022.5
$ws.Range("N74").Value = -3770.5
$ws.Range("H77").Value = 1476.9714
$ws.Range("J77").Value = 2022.5
$ws.Range("L77").Value = 10112.5
$ws.Range("N77").Value = -18848.5
$ws.Range("H102").Value = 1943
$ws.Range("I102").Value = 1368.5555
$ws.Range("K102").Value = 1368.5555
$ws.Range("M102").Value = 253.4445000000001
$ws.Range("H132").Value = 504627.7
$ws.Range("J132").Value = 7089.143
$ws.Range("L132").Value = 21267.429
$ws.Range("N132").Value = -26327.429
$ws.Range("H136").Value = 3396.742
$ws.Range("I136").Value = 2967.2
$ws.Range("K136").Value = 8901.599999999999
$ws.Range("M136").Value = -6351.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1890.92
$ws.Range("I86").Value = 1643.9474
$ws.Range("K86").Value = 1643.9474
$ws.Range("M86").Value = -520.9474
$ws.Range("H89").Value = 1890.92
$ws.Range("I89").Value = 1643.9474
$ws.Range("K89").Value = 8219.737000000001
$ws.Range("M89").Value = -2603.737000000001
$ws.Range("H134").Value = 2154045.8
$ws.Range("I134").Value = 2566912
$ws.Range("J134").Value = 7141.8
$ws.Range("K134").Value = 7700736
$ws.Range("L134").Value = 21425.4
$ws.Range("M134").Value = -7698201
$ws.Range("N134").Value = -26495.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2205.954
$ws.Range("I31").Value = 1314.122
$ws.Range("K31").Value = 1314.122
$ws.Range("M31").Value = -1019.122
$ws.Range("H34").Value = 2205.954
$ws.Range("I34").Value = 1314.122
$ws.Range("K34").Value = 1314.122
$ws.Range("M34").Value = -1112.122
$ws.Range("H99").Value = 2815.3333
$ws.Range("J99").Value = 2827
$ws.Range("L99").Value = 2827
$ws.Range("N99").Value = -5823
$ws.Range("H126").Value = 2815.3333
$ws.Range("J126").Value = 2827
$ws.Range("L126").Value = 8481
$ws.Range("N126").Value = -13421

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 904.1667
$ws.Range("I17").Value = 1025.1
$ws.Range("J17").Value = 299.5
$ws.Range("K17").Value = 3075.3
$ws.Range("L17").Value = 898.5
$ws.Range("M17").Value = -2906.3
$ws.Range("N17").Value = -1236.5
$ws.Range("H32").Value = 1679.4
$ws.Range("J32").Value = 2500
$ws.Range("L32").Value = 7500
$ws.Range("N32").Value = -8066
$ws.Range("H33").Value = 477.45456
$ws.Range("I33").Value = 42.857143
$ws.Range("J33").Value = 1238
$ws.Range("K33").Value = 257.142858
$ws.Range("L33").Value = 7428
$ws.Range("M33").Value = 25.85714200000001
$ws.Range("N33").Value = -7994
$ws.Range("H34").Value = 2695.647
$ws.Range("J34").Value = 3622.4167
$ws.Range("L34").Value = 10867.2501
$ws.Range("N34").Value = -11035.2501
$ws.Range("H38").Value = 481.91666
$ws.Range("I38").Value = 80.57143000000001
$ws.Range("J38").Value = 1043.8
$ws.Range("K38").Value = 241.71429
$ws.Range("L38").Value = 3131.4
$ws.Range("M38").Value = 105.28571
$ws.Range("N38").Value = -3825.4
$ws.Range("H39").Value = 4273.4165
$ws.Range("J39").Value = 4273.4165
$ws.Range("L39").Value = 12820.2495
$ws.Range("N39").Value = -13408.2495
$ws.Range("H46").Value = 1749.3334
$ws.Range("I46").Value = 1199.2
$ws.Range("J46").Value = 4500
$ws.Range("K46").Value = 3597.6
$ws.Range("L46").Value = 13500
$ws.Range("M46").Value = -3506.6
$ws.Range("N46").Value = -13682
$ws.Range("H55").Value = 3471.2307
$ws.Range("I55").Value = 1320.2
$ws.Range("J55").Value = 4815.625
$ws.Range("K55").Value = 3960.6
$ws.Range("L55").Value = 14446.875
$ws.Range("M55").Value = -3783.6
$ws.Range("N55").Value = -14800.875
$ws.Range("H68").Value = 2160.923
$ws.Range("J68").Value = 2249.2
$ws.Range("L68").Value = 6747.599999999999
$ws.Range("N68").Value = -8369.599999999999
$ws.Range("H71").Value = 2160.923
$ws.Range("J71").Value = 2249.2
$ws.Range("L71").Value = 20242.8
$ws.Range("N71").Value = -28354.8
$ws.Range("H97").Value = 500.66666
$ws.Range("I97").Value = 451.5
$ws.Range("J97").Value = 599
$ws.Range("K97").Value = 1354.5
$ws.Range("L97").Value = 1797
$ws.Range("M97").Value = -858.5
$ws.Range("N97").Value = -2789
$ws.Range("H117").Value = 168854.67
$ws.Range("I117").Value = 1043
$ws.Range("J117").Value = 336666.34
$ws.Range("K117").Value = 3129
$ws.Range("L117").Value = 1009999.02
$ws.Range("M117").Value = 313
$ws.Range("N117").Value = -1016883.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3100
$ws.Range("I80").Value = 3400
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3400
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2402
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 3100
$ws.Range("I83").Value = 3400
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 17000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -12008
$ws.Range("N83").Value = -24984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H59").Value = 59931.332
$ws.Range("J59").Value = 59931.332
$ws.Range("L59").Value = 59931.332
$ws.Range("N59").Value = -61239.332
$ws.Range("H136").Value = 46320
$ws.Range("I136").Value = 126959.336
$ws.Range("J136").Value = 6000.3335
$ws.Range("K136").Value = 380878.008
$ws.Range("L136").Value = 18001.0005
$ws.Range("M136").Value = -378328.008
$ws.Range("N136").Value = -23101.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6788.5713
$ws.Range("I62").Value = 7998.5
$ws.Range("J62").Value = 5881.125
$ws.Range("K62").Value = 7998.5
$ws.Range("L62").Value = 5881.125
$ws.Range("M62").Value = -7374.5
$ws.Range("N62").Value = -7129.125
$ws.Range("H65").Value = 6788.5713
$ws.Range("I65").Value = 7998.5
$ws.Range("J65").Value = 5881.125
$ws.Range("K65").Value = 39992.5
$ws.Range("L65").Value = 29405.625
$ws.Range("M65").Value = -36872.5
$ws.Range("N65").Value = -35645.625
$ws.Range("H122").Value = 50006370
$ws.Range("I122").Value = 66673092
$ws.Range("J122").Value = 6199.6
$ws.Range("K122").Value = 200019276
$ws.Range("L122").Value = 18598.8
$ws.Range("M122").Value = -200016826
$ws.Range("N122").Value = -23498.8
$ws.Range("H136").Value = 32177.229
$ws.Range("I136").Value = 3251.7144
$ws.Range("K136").Value = 9755.143199999999
$ws.Range("M136").Value = -7205.143199999999
